# Update "想去人数" (want-to-go count) figures for the "展览" and "全部类型"
# sheets to reflect newly generated output (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F2").Value = 3461   # 南宁·第二届北极光动漫展: 3457 -> 3461
    $ws.Range("F4").Value = 73     # 南宁·原神x星铁x绝区零同人ONLY3.0: 72 -> 73
    $ws.Range("F5").Value = 1893   # 南宁·2024良牙动漫秋季盛典（秋典）: 1874 -> 1893
    $ws.Range("F6").Value = 131    # 南宁·快看漫画动漫游戏嘉年华 KKWORLD-mini: 128 -> 131
}

# "展览" sheet: 南宁·万圣漫控嘉年华10 is row 7
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F7").Value = 346   # 344 -> 346

# "全部类型" sheet: 南宁·万圣漫控嘉年华10 is row 8
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F8").Value = 346    # 344 -> 346
